$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 43
$ws1.Range("F4").Value = 2131
$ws1.Range("F5").Value = 173
$ws1.Range("F6").Value = 361

# Sheet "全部类型" (all categories)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 43
$ws4.Range("F4").Value = 2131
$ws4.Range("F5").Value = 173
$ws4.Range("F7").Value = 361
